$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.504.32"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.626.14"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'213.72"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'19.22"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").Value = "'0.0853"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.854.56"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "1.621.28"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'4.05"
$ws.Range("D15").Value = "'0.512"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "'64.01"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "'235.29"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "26.509.93"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'7.77"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "0.0₃0727"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'146.82"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").Value = "'15.65"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "'0.0498"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").Value = "1.519.12"
$ws.Range("E32").Value = "  +5.16%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "'0.836"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").Value = "'5.87"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "1.765.42"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").Value = "'0.760"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "'0.909"
$ws.Range("D47").Value = "'89.94"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = "  -0.58%  "
